$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(5).Copy()
$ws.Columns.Item(4).PasteSpecial(-4122)
Write-Host "D7:" $ws.Range("D7").NumberFormat
